# Update evaluation metrics across the three worksheets to reflect the
# refactored training/evaluation results.

$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5767790262172284
$ws1.Range("C2").Value = 0.6846846846846847
$ws1.Range("D2").Value = 0.2846441947565543
$ws1.Range("E2").Value = 0.4021164021164021
$ws1.Range("F2").Value = 0.5767790262172284
$ws1.Range("G2").Value = 152
$ws1.Range("H2").Value = 70
$ws1.Range("I2").Value = 464
$ws1.Range("J2").Value = 382

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2: class "0.0"
$ws2.Range("B2").Value = 0.5484633569739953
$ws2.Range("C2").Value = 0.8689138576779026
$ws2.Range("D2").Value = 0.672463768115942
$ws2.Range("E2").Value = 534

# row 3: class "1.0"
$ws2.Range("B3").Value = 0.6846846846846847
$ws2.Range("C3").Value = 0.2846441947565543
$ws2.Range("D3").Value = 0.4021164021164021
$ws2.Range("E3").Value = 534

# row 4: accuracy
$ws2.Range("B4").Value = 0.5767790262172284
$ws2.Range("C4").Value = 0.5767790262172284
$ws2.Range("D4").Value = 0.5767790262172284
$ws2.Range("E4").Value = 0.5767790262172284

# row 5: macro avg
$ws2.Range("B5").Value = 0.61657402082934
$ws2.Range("C5").Value = 0.5767790262172284
$ws2.Range("D5").Value = 0.537290085116172
$ws2.Range("E5").Value = 1068

# row 6: weighted avg
$ws2.Range("B6").Value = 0.61657402082934
$ws2.Range("C6").Value = 0.5767790262172284
$ws2.Range("D6").Value = 0.537290085116172
$ws2.Range("E6").Value = 1068

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2: Actual 0
$ws3.Range("B2").Value = 464
$ws3.Range("C2").Value = 70

# row 3: Actual 1
$ws3.Range("B3").Value = 382
$ws3.Range("C3").Value = 152
